$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD/AE/AF = Wins/Losses/Ties ---
# Copy the formatting of the existing header cell (AC1, style index 1:
# bold font, thin border, centered/top-aligned) onto the three new
# header cells before writing their text, so the new header cells match
# the look of the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-48): season record for every player row ---
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 86
    $ws.Cells.Item($r, 31).Value = 75
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Host "done"
